$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula = "'25.964.37"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Formula = "'  +1.56%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Formula = "'1.596.88"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Formula = "'  +1.44%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Formula = "'  +0.50%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Formula = "'210.88"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Formula = "'  +1.08%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("E6").Formula = "'  +0.50%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Formula = "'0.483"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Formula = "'  +0.12%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Formula = "'  +0.13%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Formula = "'  -1.09%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Formula = "'18.02"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Formula = "'  -0.50%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Formula = "'  +3.41%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Formula = "'1.820.85"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Formula = "'  +1.69%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Formula = "'1.595.69"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Formula = "'  +1.42%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Formula = "'4.00"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Formula = "'  -0.41%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Formula = "'0.513"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Formula = "'  +0.07%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Formula = "'25.978.21"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Formula = "'  +1.82%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Formula = "'60.04"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Formula = "'  +0.18%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Formula = "'0.0₃0721"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Formula = "'  +0.03%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Formula = "'  +0.13%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Formula = "'200.62"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Formula = "'  +5.72%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Formula = "'  +1.36%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Formula = "'  -1.23%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Formula = "'  +1.35%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Formula = "'1.78"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Formula = "'  +6.37%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Formula = "'141.96"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Formula = "'  +0.70%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Formula = "'  +0.49%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Formula = "'  -8.58%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Formula = "'15.09"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Formula = "'  -0.17%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Formula = "'6.45"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Formula = "'  +0.20%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Formula = "'  +0.66%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Formula = "'0.0475"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Formula = "'  +1.06%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Formula = "'  +0.24%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Formula = "'2.95"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Formula = "'  -2.02%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Formula = "'  -1.06%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Formula = "'  +2.84%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Formula = "'1.125.10"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Formula = "'  +2.82%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Formula = "'0.0162"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Formula = "'  +8.18%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Formula = "'  +0.41%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Formula = "'  -0.91%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Formula = "'0.786"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Formula = "'  +0.39%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Formula = "'  -2.53%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Formula = "'0.782"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Formula = "'  -2.66%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Formula = "'1.732.16"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Formula = "'  +1.93%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Formula = "'  -0.58%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Formula = "'92.63"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Formula = "'  -0.88%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Formula = "'  -0.18%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Formula = "'53.38"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Formula = "'  +0.73%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Formula = "'  -0.72%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Formula = "'  +0.57%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Formula = "'  +0.65%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Formula = "'7.15"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Formula = "'  -2.41%  "
$ws.Range("E51").Style = "Normal"
